$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The label column (B) was relabelled: rows 2-216 flip 0 -> 1,
# and rows 217-431 flip 1 -> 0.
$ws.Range("B2:B216").Value = 1
$ws.Range("B217:B431").Value = 0

# Update the active selection to match the saved state (E16).
$ws.Range("E16").Select() | Out-Null
